$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for each data row (2-8).
# The serial value 46074 should be bumped to 46075 (one day later) for all of them.
foreach ($row in 2..8) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46074) {
        $cell.Value2 = 46075
    }
}
